# TODO list update: mark "DEATH" (fireball/bat kill logic) and
# "set up triggers for room changes" as done by Tom, per the commit message
# ("trigger problem fixed ... made it so fireball will be destroyed but not
# kill the bats if they aren't tired ... works now, is in player lose hp").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: "set up triggers for room changes" -> completed by tom, status "done"
$ws.Range("B13").Value = "tom"
$ws.Range("D13").Value = "done"

# Row 9: "DEATH" -> completed by tom, with a status note
$ws.Range("B9").Value = "tom"
$ws.Range("D9").Value = "works now, is in player lose hp "

# Move the active selection, matching where the author left off editing
[void]$ws.Range("B14").Select()
